# Add "Credit On Tenure Basis" column to the "All_Without_Probation" sheet.
# A new column is inserted before column V; every data row (2-83) gets
# "Yes" in that column, and the header (row 1) gets the new column title.
# This mirrors the already-existing layout on the "All_Scenarios" sheet,
# which has the same "Credit On Tenure Basis" / "Credit From Year" /
# "Credit To Year" / "Credit No of Leaves" columns (V:Y).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All_Without_Probation")

# Insert a new column at V; existing V:X ("Credit From Year", "Credit To
# Year", "Credit No of Leaves") shift right to W:Y.
$ws.Columns("V").Insert()

# New header cell for the inserted column.
$ws.Range("V1").Value = "Credit On Tenure Basis"

# Every scenario row gets "Yes" for the new "on tenure basis" flag.
$ws.Range("V2:V83").Value = "Yes"
$ws.Range("V2:V83").WrapText = $false

# Re-wrapping column B (description text) against the (slightly) changed
# layout pushes a handful of rows to wrap one extra line; bump those rows'
# heights to match.
$ws.Rows.Item(4).RowHeight = 113.95
$ws.Rows.Item(10).RowHeight = 102.7
$ws.Rows.Item(15).RowHeight = 147.7
$ws.Rows.Item(18).RowHeight = 136.45
$ws.Rows.Item(19).RowHeight = 136.45
$ws.Rows.Item(20).RowHeight = 147.7
$ws.Rows.Item(21).RowHeight = 136.45
$ws.Rows.Item(22).RowHeight = 125.2
$ws.Rows.Item(25).RowHeight = 136.45
$ws.Rows.Item(30).RowHeight = 147.7
$ws.Rows.Item(35).RowHeight = 147.7
$ws.Rows.Item(40).RowHeight = 147.7
$ws.Rows.Item(45).RowHeight = 113.95
$ws.Rows.Item(51).RowHeight = 102.7
$ws.Rows.Item(56).RowHeight = 147.7
$ws.Rows.Item(59).RowHeight = 136.45
$ws.Rows.Item(60).RowHeight = 136.45
$ws.Rows.Item(61).RowHeight = 147.7
$ws.Rows.Item(62).RowHeight = 136.45
$ws.Rows.Item(63).RowHeight = 125.2
$ws.Rows.Item(66).RowHeight = 136.45
$ws.Rows.Item(71).RowHeight = 147.7
$ws.Rows.Item(76).RowHeight = 147.7
$ws.Rows.Item(81).RowHeight = 147.7

# Leave the new column selected, matching the end state after the edit.
$ws.Range("V2:V83").Select() | Out-Null
